$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Training RMSE block (rows 5-7) ---
$ws.Range("C5").Value = 0.819402704056875
$ws.Range("D5").Value = 0.775051439637763
$ws.Range("E5").Value = 0.801025969828675
$ws.Range("F5").Value = 0.928403411130697

$ws.Range("C6").Value = 0.562082010351434
$ws.Range("D6").Value = 0.552962097649074
$ws.Range("E6").Value = 0.585287783208194
$ws.Range("F6").Value = 0.854073613966128

$ws.Range("C7").Value = 0.222858457658676
$ws.Range("D7").Value = 0.211689990670374
$ws.Range("E7").Value = 0.298868873298056
$ws.Range("F7").Value = 0.874983841162522

# --- Update Validation RMSE block (rows 12-14) ---
$ws.Range("C12").Value = 0.958707283679204
$ws.Range("D12").Value = 0.960188377135612
$ws.Range("E12").Value = 0.949663040567693
$ws.Range("F12").Value = 0.983660028016208

$ws.Range("C13").Value = 0.865610485637818
$ws.Range("D13").Value = 0.867606451132059
$ws.Range("E13").Value = 0.87093870187506
$ws.Range("F13").Value = 0.95669128535637

$ws.Range("C14").Value = 0.705492603737855
$ws.Range("D14").Value = 0.694408799265067
$ws.Range("E14").Value = 0.702073984969587
$ws.Range("F14").Value = 0.945423353622081

# --- Fill in Testing RMSE block (rows 19-21), previously empty ---
$ws.Range("C19").Value = 1.08157893901907
$ws.Range("D19").Value = 1.08538228022025
$ws.Range("E19").Value = 1.07293191044207
$ws.Range("F19").Value = 1.11604897913846

$ws.Range("C20").Value = 0.962551200148463
$ws.Range("D20").Value = 0.964462006793618
$ws.Range("E20").Value = 0.97463378415388
$ws.Range("F20").Value = 1.08147900935051

$ws.Range("C21").Value = 0.763065037439799
$ws.Range("D21").Value = 0.775224764725676
$ws.Range("E21").Value = 0.786992541953138
$ws.Range("F21").Value = 1.0698128993807

# --- Slightly narrower columns A and B (minor manual resize in the original edit) ---
$ws.Columns.Item(1).ColumnWidth = 7.333333333333334
$ws.Columns.Item(2).ColumnWidth = 18.0

# --- Move the active selection to H18 ---
$ws.Range("H18").Select()
